# record add one hour
# Append a new time-log entry (row 52) below the existing table: a new
# date/start/finish row whose duration (Finish-Start) works out to one
# extra hour logged, consistent with the "record add one hour" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the formatting of the last existing row (51) down onto the new
# row (52) first, so the new cells pick up the same date/time number
# formats and fonts already used throughout the table.
$ws.Range("A51:D51").Copy() | Out-Null
$ws.Range("A52:D52").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Rows.Item(52).RowHeight = 15.75

# New entry: 2023-10-25 (serial 45224), 21:00 -> 22:00 (one hour).
$ws.Range("A52").NumberFormat = "mm-dd-yy"
$ws.Range("A52").Value = 45224
$ws.Range("B52").Value = 0.875
$ws.Range("C52").Value = 0.91666666666666663

# Duration formula, matching the pattern used by the rest of column D.
$ws.Range("D52").Formula = "=C52-B52"

$excel.CutCopyMode = 0

# Match the saved file's active selection on the freshly-added row.
$null = $ws.Range("E52").Select()
